$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.433.37"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.27"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.20"
$ws.Range("E5").Value = "  -4.06%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  -3.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06797"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.72"
$ws.Range("E10").Value = "  -6.10%  "
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07741"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.847.53"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.011"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  -3.62%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007974"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.468.24"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.075.95"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.622"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.580"
$ws.Range("E23").Value = "  -4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.987"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.55"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.176"
$ws.Range("E26").Value = "  -8.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.654"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.01"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.67"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.162"
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.126"
$ws.Range("E31").Value = "  -4.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08710"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04829"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.130"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7186"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.842"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.106"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.226"
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4860"
$ws.Range("E40").Value = "  -4.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9132"
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.94"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.721"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4173"
$ws.Range("E46").Value = "  -5.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05923"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.119"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1236"
$ws.Range("E49").Value = "  -7.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.05"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8859"
$ws.Range("E51").Value = "  +0.72%  "
